# Insert a new data row at row 307, shifting existing rows 307:380 down to 308:381.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(307).Insert()

$ws.Cells.Item(307, 1).Value = 4
$ws.Cells.Item(307, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(307, 3).Value = "Los Lagos"
$ws.Cells.Item(307, 4).Value = 44782
$ws.Cells.Item(307, 5).Value = 10
$ws.Cells.Item(307, 6).Value = 100112008
$ws.Cells.Item(307, 7).Value = "Coliflor"
$ws.Cells.Item(307, 8).Value = "Sin especificar"
$ws.Cells.Item(307, 9).Value = "Primera"
$ws.Cells.Item(307, 10).Value = 1400
$ws.Cells.Item(307, 11).Value = 1500
$ws.Cells.Item(307, 12).Value = 1700
$ws.Cells.Item(307, 13).Value = 1600
$ws.Cells.Item(307, 14).Value = "$/unidad"
$ws.Cells.Item(307, 15).Value = "Región Metropolitana"
$ws.Cells.Item(307, 16).Value = 1600
$ws.Cells.Item(307, 17).Value = 1
$ws.Cells.Item(307, 18).Value = "Hortaliza"
